$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 10
$ws.Range("C2").Value = 20
$ws.Range("D2").Value = 30
$ws.Range("E2").Value = 40

$ws.Range("R15").Font.Underline = $true
$ws.Range("R15").Select()
